$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New activity-log rows appended to the bottom of the log table (rows 30-35),
# matching the existing table's layout: User | Message | Date | Time
$rows = @(
    @{ User = "seyntt1"; Message = "seyntt1 logged in"; Date = 45789; Time = 45789.8446643518 },
    @{ User = "seyntt1"; Message = "seyntt1 logged in"; Date = 45789; Time = 45789.8475578704 },
    @{ User = "seyntt1"; Message = "seyntt1 logged in"; Date = 45789; Time = 45789.8509259259 },
    @{ User = "seyntt1"; Message = "seyntt1 logged in"; Date = 45789; Time = 45789.8518865741 },
    @{ User = "we123";   Message = "we123 logged in";   Date = 45789; Time = 45789.8534027778 },
    @{ User = "seyntt1"; Message = "seyntt1 logged in"; Date = 45789; Time = 45789.8743865741 }
)

$lastRow = 29
$templateFormat = $ws.Range("C" + $lastRow + ":D" + $lastRow).NumberFormat

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $entry = $rows[$i]

    $ws.Range("A$r").Value = $entry.User
    $ws.Range("B$r").Value = $entry.Message
    $ws.Range("C$r").Value = $entry.Date
    $ws.Range("D$r").Value = $entry.Time

    $ws.Range("C$r`:D$r").NumberFormat = $templateFormat
}
